$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2863.3
$ws.Range("C2").Value = 2805.15
$ws.Range("D2").Value = 2849.05
$ws.Range("E2").Value = 2848.95
$ws.Range("F2").Value = 15
$ws.Range("G2").Value = 2810

$ws.Range("B3").Value = 455.15
$ws.Range("C3").Value = 443.3
$ws.Range("D3").Value = 453.95
$ws.Range("E3").Value = 454.05
$ws.Range("F3").Value = 32
$ws.Range("G3").Value = 444.75

$ws.Range("B4").Value = 1692.7
$ws.Range("C4").Value = 1678.4
$ws.Range("D4").Value = 1688
$ws.Range("E4").Value = 1685.8
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 1679.5

$ws.Range("B5").Value = 7343
$ws.Range("C5").Value = 7215.45
$ws.Range("D5").Value = 7315
$ws.Range("E5").Value = 7327.75
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 7221.85

$ws.Range("B6").Value = 243
$ws.Range("C6").Value = 239.4
$ws.Range("D6").Value = 241.1
$ws.Range("E6").Value = 241.4
$ws.Range("F6").Value = 45
$ws.Range("G6").Value = 239.8

$ws.Range("B7").Value = 232.75
$ws.Range("C7").Value = 229.35
$ws.Range("D7").Value = 230.85
$ws.Range("E7").Value = 231.1
$ws.Range("F7").Value = 110
$ws.Range("G7").Value = 232.55

$ws.Range("B8").Value = 379.9
$ws.Range("C8").Value = 374.25
$ws.Range("D8").Value = 376.2
$ws.Range("E8").Value = 376
$ws.Range("F8").Value = 88
$ws.Range("G8").Value = 377.8

$ws.Range("B9").Value = 730
$ws.Range("C9").Value = 721.8
$ws.Range("D9").Value = 724
$ws.Range("E9").Value = 726.4
$ws.Range("F9").Value = 19
$ws.Range("G9").Value = 726.9

$ws.Range("B10").Value = 4178
$ws.Range("C10").Value = 4096.05
$ws.Range("D10").Value = 4132
$ws.Range("E10").Value = 4143.5
$ws.Range("F10").Value = 6
$ws.Range("G10").Value = 4116.4

$ws.Range("B11").Value = 157.75
$ws.Range("C11").Value = 154.9
$ws.Range("D11").Value = 156.45
$ws.Range("E11").Value = 156.15
$ws.Range("F11").Value = 130
$ws.Range("G11").Value = 157.6

$ws.Range("B12").Value = 1481.85
$ws.Range("C12").Value = 1457.55
$ws.Range("D12").Value = 1468
$ws.Range("E12").Value = 1466.1
$ws.Range("F12").Value = 14
$ws.Range("G12").Value = 1470.05

$ws.Range("B13").Value = 1714.9
$ws.Range("C13").Value = 1700.3
$ws.Range("D13").Value = 1708.65
$ws.Range("E13").Value = 1709.25
$ws.Range("F13").Value = 125
$ws.Range("G13").Value = 1702.5

$ws.Range("B14").Value = 618
$ws.Range("C14").Value = 608.1
$ws.Range("D14").Value = 614.25
$ws.Range("E14").Value = 614.85
$ws.Range("F14").Value = 39
$ws.Range("G14").Value = 611

$ws.Range("B15").Value = 1002.25
$ws.Range("C15").Value = 993.7
$ws.Range("D15").Value = 998.4
$ws.Range("E15").Value = 996.6
$ws.Range("F15").Value = 122
$ws.Range("G15").Value = 1000.5

$ws.Range("B16").Value = 1606.4
$ws.Range("C16").Value = 1587.05
$ws.Range("D16").Value = 1595.9
$ws.Range("E16").Value = 1598.95
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 1604

$ws.Range("B17").Value = 1554.9
$ws.Range("C17").Value = 1540
$ws.Range("D17").Value = 1542.65
$ws.Range("E17").Value = 1542.9
$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 1550.7

$ws.Range("B18").Value = 750
$ws.Range("C18").Value = 735.9
$ws.Range("D18").Value = 747
$ws.Range("E18").Value = 748.1
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 737.7

$ws.Range("B19").Value = 540
$ws.Range("C19").Value = 525.65
$ws.Range("D19").Value = 536.5
$ws.Range("E19").Value = 536.1
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 526.35

$ws.Range("B20").Value = 1758
$ws.Range("C20").Value = 1715
$ws.Range("D20").Value = 1725
$ws.Range("E20").Value = 1729.4
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 1727.1

$ws.Range("B21").Value = 281.55
$ws.Range("C21").Value = 274.5
$ws.Range("D21").Value = 277
$ws.Range("E21").Value = 276.7
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 275.05

$ws.Range("B22").Value = 313.1
$ws.Range("C22").Value = 308.15
$ws.Range("D22").Value = 310.5
$ws.Range("E22").Value = 311.15
$ws.Range("F22").Value = 126
$ws.Range("G22").Value = 312.9

$ws.Range("B23").Value = 2597.8
$ws.Range("C23").Value = 2579.3
$ws.Range("D23").Value = 2587
$ws.Range("E23").Value = 2584.95
$ws.Range("F23").Value = 54
$ws.Range("G23").Value = 2585.75

$ws.Range("B24").Value = 649.6
$ws.Range("C24").Value = 639.6
$ws.Range("D24").Value = 641.7
$ws.Range("E24").Value = 642.05
$ws.Range("F24").Value = 132
$ws.Range("G24").Value = 647.3

$ws.Range("B25").Value = 714.5
$ws.Range("C25").Value = 706.85
$ws.Range("D25").Value = 711.9
$ws.Range("E25").Value = 712.2
$ws.Range("F25").Value = 3
$ws.Range("G25").Value = 711.15

$ws.Range("B26").Value = 1123
$ws.Range("C26").Value = 1089
$ws.Range("D26").Value = 1107.85
$ws.Range("E26").Value = 1104.1
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 1092.7

$ws.Range("B27").Value = 802.9
$ws.Range("C27").Value = 762.5
$ws.Range("D27").Value = 779.4
$ws.Range("E27").Value = 779.95
$ws.Range("F27").Value = 411
$ws.Range("G27").Value = 764.5

$ws.Range("B28").Value = 336
$ws.Range("C28").Value = 326.05
$ws.Range("D28").Value = 331.5
$ws.Range("E28").Value = 332.15
$ws.Range("F28").Value = 199
$ws.Range("G28").Value = 328.05

$ws.Range("B29").Value = 141.25
$ws.Range("C29").Value = 137.2
$ws.Range("D29").Value = 139.3
$ws.Range("E29").Value = 139.6
$ws.Range("F29").Value = 492
$ws.Range("G29").Value = 137.55

$ws.Range("B30").Value = 10526
$ws.Range("C30").Value = 10361.55
$ws.Range("D30").Value = 10482.35
$ws.Range("E30").Value = 10503.05
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 10395.4

$ws.Range("D2").Select()
